$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" timestamp
# shared by the 02993fd4... and 8f234ff7... records, which were generated in
# the same handback run (Overview rows 2 & 4, de-de H column rows 2 & 4).
$wsOverview.Range("G2").Value = "2016-09-06 00:18:28"
$wsOverview.Range("G4").Value = "2016-09-06 00:18:28"
$wsDeDe.Range("H2").Value = "2016-09-06 00:18:28"
$wsDeDe.Range("H4").Value = "2016-09-06 00:18:28"

# Status changed from handoff type "ht" to "mt" for both records, in both
# language sheets.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# zh-cn Correspond Handoff / Handback datetimes
$wsZhCn.Range("H2").Value = "2016-09-06 00:18:24"
$wsZhCn.Range("H4").Value = "2016-09-06 00:18:24"
$wsZhCn.Range("K2").Value = "2016-09-06 00:18:41"
$wsZhCn.Range("K4").Value = "2016-09-06 00:18:41"

# de-de Correspond Handback datetime
$wsDeDe.Range("K2").Value = "2016-09-06 00:18:49"
$wsDeDe.Range("K4").Value = "2016-09-06 00:18:49"
